# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to the freshly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 463
$ws1.Range("F4").Value  = 7746
$ws1.Range("F5").Value  = 91
$ws1.Range("F8").Value  = 26
$ws1.Range("F10").Value = 445
$ws1.Range("F14").Value = 63
$ws1.Range("F15").Value = 65
$ws1.Range("F17").Value = 5607
$ws1.Range("F18").Value = 155
$ws1.Range("F19").Value = 218
$ws1.Range("F20").Value = 1014
$ws1.Range("F22").Value = 326

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 463
$ws4.Range("F4").Value  = 7746
$ws4.Range("F5").Value  = 91
$ws4.Range("F8").Value  = 26
$ws4.Range("F10").Value = 445
$ws4.Range("F14").Value = 63
$ws4.Range("F15").Value = 65
$ws4.Range("F18").Value = 5607
$ws4.Range("F20").Value = 155
$ws4.Range("F21").Value = 218
$ws4.Range("F22").Value = 1014
$ws4.Range("F24").Value = 326
